$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting (it is displayed
# as plain text in the source data, e.g. "1.000", "28.429.14"), so force the
# column to Text format before writing values to avoid Excel re-interpreting
# these strings as numbers/dates and dropping trailing zeros, etc.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.429.14'
$ws.Range('E2').Value = '  +4.16%  '
$ws.Range('D3').Value = '1.807.52'
$ws.Range('E3').Value = '  +1.61%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '316.38'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '0.5462'
$ws.Range('E7').Value = '  +4.85%  '
$ws.Range('D8').Value = '0.3850'
$ws.Range('E8').Value = '  +6.65%  '
$ws.Range('D9').Value = '0.07600'
$ws.Range('E9').Value = '  +2.89%  '
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('D11').Value = '1.127'
$ws.Range('E11').Value = '  +3.01%  '
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('E13').Value = '  +2.85%  '
$ws.Range('E14').Value = '  +2.10%  '
$ws.Range('D15').Value = '7.385'
$ws.Range('E15').Value = '  +5.55%  '
$ws.Range('D16').Value = '1.803.56'
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('D17').Value = '92.23'
$ws.Range('E17').Value = '  +4.20%  '
$ws.Range('E18').Value = '  +2.28%  '
$ws.Range('D19').Value = '0.06445'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').Value = '0.9994'
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').Value = '17.39'
$ws.Range('E21').Value = '  +4.00%  '
$ws.Range('D22').Value = '5.987'
$ws.Range('E22').Value = '  +2.18%  '
$ws.Range('D23').Value = '28.430.76'
$ws.Range('E23').Value = '  +3.79%  '
$ws.Range('D24').Value = '11.45'
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('D25').Value = '2.135'
$ws.Range('E25').Value = '  +3.28%  '
$ws.Range('D26').Value = '158.72'
$ws.Range('E26').Value = '  +3.29%  '
$ws.Range('D27').Value = '20.67'
$ws.Range('E27').Value = '  +2.76%  '
$ws.Range('D28').Value = '2.409'
$ws.Range('E28').Value = '  +2.47%  '
$ws.Range('D29').Value = '2.012.02'
$ws.Range('E29').Value = '  +1.12%  '
$ws.Range('D30').Value = '123.81'
$ws.Range('E30').Value = '  +1.79%  '
$ws.Range('D31').Value = '1.126'
$ws.Range('E31').Value = '  +6.04%  '
$ws.Range('D32').Value = '0.1019'
$ws.Range('E32').Value = '  +4.52%  '
$ws.Range('D33').Value = '5.758'
$ws.Range('E33').Value = '  +3.68%  '
$ws.Range('D34').Value = '3.694'
$ws.Range('E34').Value = '  +2.70%  '
$ws.Range('D35').Value = '0.2326'
$ws.Range('E35').Value = '  +14.74%  '
$ws.Range('D36').Value = '0.06387'
$ws.Range('E36').Value = '  +6.54%  '
$ws.Range('D37').Value = '0.02329'
$ws.Range('E37').Value = '  +3.95%  '
$ws.Range('D38').Value = '8.841'
$ws.Range('E38').Value = '  +9.79%  '
$ws.Range('D39').Value = '5.100'
$ws.Range('E39').Value = '  +5.07%  '
$ws.Range('E40').Value = '  +3.74%  '
$ws.Range('D41').Value = '0.6429'
$ws.Range('E41').Value = '  +4.41%  '
$ws.Range('B42').Value = 'Frax'
$ws.Range('C42').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D42').Value = '0.9992'
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '1.161'
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('D44').Value = '1.384'
$ws.Range('E44').Value = '  -3.14%  '
$ws.Range('E45').Value = '  +2.20%  '
$ws.Range('D46').Value = '0.5982'
$ws.Range('E46').Value = '  +3.64%  '
$ws.Range('E47').Value = '  +1.47%  '
$ws.Range('D48').Value = '126.49'
$ws.Range('E48').Value = '  +3.88%  '
$ws.Range('D49').Value = '1.991'
$ws.Range('E49').Value = '  +5.25%  '
$ws.Range('E50').Value = '  +3.42%  '
$ws.Range('D51').Value = '0.06908'
$ws.Range('E51').Value = '  +2.84%  '
